# Apply the cryptos list refresh (prices / 1h volume %, and the swapped
# dogwifhat / Hedera rows) as captured by the GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column "Price" cells whose new value is purely numeric-looking would
# otherwise be auto-coerced to a Number by Excel (and lose e.g. a trailing
# zero, "352.10" -> 352.1). Force those specific cells to Text first so the
# stored value matches the source string exactly, just like the others that
# already read as Text because of repeated "." separators.
$textPriceCells = @(5, 6, 7, 10, 11, 14, 18, 19, 20, 21, 23, 24, 25, 26, 28, 29, 30, 33, 34, 35, 37, 39, 40, 41, 43, 44, 45, 46)
foreach ($r in $textPriceCells) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "65.896.90"
$ws.Range("E2").Value = "  +0.34%  "

$ws.Range("D3").Value = "2.667.57"
$ws.Range("E3").Value = "  -0.42%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "598.72"

$ws.Range("D6").Value = "157.84"
$ws.Range("E6").Value = "  +0.63%  "

$ws.Range("D7").Value = "0.652"
$ws.Range("E7").Value = "  +4.33%  "

$ws.Range("E9").Value = "  -2.03%  "

$ws.Range("D10").Value = "0.402"
$ws.Range("E10").Value = "  +0.31%  "

$ws.Range("D11").Value = "5.87"
$ws.Range("E11").Value = "  +0.17%  "

$ws.Range("E12").Value = "  +1.64%  "

$ws.Range("E13").Value = "  -0.93%  "

$ws.Range("D14").Value = "0.0000195"
$ws.Range("E14").Value = "  -1.69%  "

$ws.Range("D15").Value = "3.146.50"
$ws.Range("E15").Value = "  -0.31%  "

$ws.Range("D16").Value = "65.801.24"
$ws.Range("E16").Value = "  +0.43%  "

$ws.Range("D17").Value = "2.686.04"
$ws.Range("E17").Value = "  +0.27%  "

$ws.Range("D18").Value = "12.67"
$ws.Range("E18").Value = "  -1.10%  "

$ws.Range("D19").Value = "4.81"
$ws.Range("E19").Value = "  +0.07%  "

$ws.Range("D20").Value = "352.10"
$ws.Range("E20").Value = "  +0.12%  "

$ws.Range("D21").Value = "7.49"
$ws.Range("E21").Value = "  -1.19%  "

$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").Value = "69.90"
$ws.Range("E23").Value = "  +0.43%  "

$ws.Range("D24").Value = "1.83"
$ws.Range("E24").Value = "  +11.57%  "

$ws.Range("D25").Value = "0.0000112"
$ws.Range("E25").Value = "  +0.64%  "

$ws.Range("D26").Value = "9.69"
$ws.Range("E26").Value = "  +0.21%  "

$ws.Range("E27").Value = "  +2.15%  "

$ws.Range("D28").Value = "576.34"
$ws.Range("E28").Value = "  +8.48%  "

$ws.Range("D29").Value = "8.25"
$ws.Range("E29").Value = "  +1.75%  "

$ws.Range("D30").Value = "0.164"
$ws.Range("E30").Value = "  -1.86%  "

$ws.Range("E31").Value = "  -0.08%  "

$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("D33").Value = "1.83"
$ws.Range("E33").Value = "  +4.05%  "

$ws.Range("D34").Value = "6.70"
$ws.Range("E34").Value = "  +3.65%  "

$ws.Range("D35").Value = "5.59"
$ws.Range("E35").Value = "  +1.72%  "

$ws.Range("D37").Value = "20.59"
$ws.Range("E37").Value = "  +0.37%  "

$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("D39").Value = "1.96"
$ws.Range("E39").Value = "  +0.75%  "

$ws.Range("D40").Value = "154.11"
$ws.Range("E40").Value = "  -2.52%  "

$ws.Range("D41").Value = "161.45"
$ws.Range("E41").Value = "  -1.89%  "

$ws.Range("E42").Value = "  -1.21%  "

$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").Value = "0.0619"
$ws.Range("E43").Value = "  +1.68%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "2.33"
$ws.Range("E44").Value = "  +0.38%  "

$ws.Range("D45").Value = "23.09"
$ws.Range("E45").Value = "  +1.06%  "

$ws.Range("D46").Value = "0.645"
$ws.Range("E46").Value = "  +0.29%  "

$ws.Range("E47").Value = "  -1.21%  "

$ws.Range("E48").Value = "  +1.16%  "

$ws.Range("E49").Value = "  -1.49%  "

$ws.Range("D50").Value = "0.0₆0245"
$ws.Range("E50").Value = "  -6.40%  "

$ws.Range("E51").Value = "  -0.65%  "
